$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header column I2: a new repeating-schedule session ---
$ws.Range("I2").Value = "check lịch lặp 1 (T5, 30.11.23)"

# --- Consolidate the per-day attendance marks into a single column C ---
# (the "move class course to lms file" commit folds the separate weekday
#  columns D/E back into column C for rows 4,6,7,8,11; and marks rows 9/10
#  as attended in column C as well)
$ws.Range("C4").Value = 1
$ws.Range("D4").ClearContents()

$ws.Range("C6").Value = 1
$ws.Range("E6").ClearContents()

$ws.Range("C7").Value = 1
$ws.Range("E7").ClearContents()

$ws.Range("D8").ClearContents()

$ws.Range("C9").Value = 1

$ws.Range("C10").Value = 1

$ws.Range("C11").Value = 1
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()

# --- Column B width tweak ---
$ws.Columns.Item(2).ColumnWidth = 31.67

# --- Selection moves to C3:C12 ---
[void]$ws.Range("C3:C12").Select()
